# Apply cryptos list update (Mon Mar 27 23:49:36 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.167.58"
$ws.Range('D3').Value = "'1.714.81"
$ws.Range('E3').Value = '  -3.68%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'310.81"
$ws.Range('E5').Value = '  -5.67%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  +6.34%  '
$ws.Range('D8').Value = "'0.3455"
$ws.Range('E8').Value = '  -3.25%  '
$ws.Range('D9').Value = "'42.54"
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').Value = "'0.07285"
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('E11').Value = '  -5.62%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = "'19.90"
$ws.Range('E13').Value = '  -5.44%  '
$ws.Range('D14').Value = "'5.870"
$ws.Range('E14').Value = '  -3.34%  '
$ws.Range('D15').Value = "'1.713.11"
$ws.Range('E15').Value = '  -3.66%  '
$ws.Range('D16').Value = "'6.878"
$ws.Range('E16').Value = '  -5.55%  '
$ws.Range('D17').Value = "'89.12"
$ws.Range('E17').Value = '  -4.98%  '
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').Value = "'0.06360"
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  -4.00%  '
$ws.Range('D22').Value = "'5.632"
$ws.Range('E22').Value = '  -3.13%  '
$ws.Range('D23').Value = "'27.206.20"
$ws.Range('E23').Value = '  -3.10%  '
$ws.Range('E24').Value = '  -4.68%  '
$ws.Range('D25').Value = "'2.089"
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('D26').Value = "'151.72"
$ws.Range('E26').Value = '  -6.29%  '
$ws.Range('E27').Value = '  -3.64%  '
$ws.Range('D28').Value = "'1.910.35"
$ws.Range('E28').Value = '  -3.66%  '
$ws.Range('E29').Value = '  -3.09%  '
$ws.Range('D30').Value = "'120.11"
$ws.Range('E30').Value = '  -4.10%  '
$ws.Range('E31').Value = '  -8.44%  '
$ws.Range('D32').Value = "'0.09240"
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').Value = "'3.582"
$ws.Range('E33').Value = '  -2.84%  '
$ws.Range('D34').Value = "'5.341"
$ws.Range('E34').Value = '  -6.64%  '
$ws.Range('D35').Value = "'0.02203"
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('D36').Value = "'0.05925"
$ws.Range('E36').Value = '  -4.87%  '
$ws.Range('D37').Value = "'11.11"
$ws.Range('E37').Value = '  -6.77%  '
$ws.Range('D38').Value = "'0.2005"
$ws.Range('E38').Value = '  -5.29%  '
$ws.Range('D39').Value = "'1.422"
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('D40').Value = "'4.760"
$ws.Range('E40').Value = '  -5.18%  '
$ws.Range('D41').Value = "'1.0000"
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  -6.39%  '
$ws.Range('E43').Value = '  -7.40%  '
$ws.Range('D44').Value = "'7.488"
$ws.Range('E44').Value = '  -5.66%  '
$ws.Range('D45').Value = "'12.75"
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('D46').Value = "'3.584"
$ws.Range('E46').Value = '  -4.60%  '
$ws.Range('E47').Value = '  -5.21%  '
$ws.Range('D48').Value = "'118.71"
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('E49').Value = '  -6.38%  '
$ws.Range('D50').Value = "'0.06645"
$ws.Range('E51').Value = '  -5.23%  '
